$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new columns of header values (P1=14, Q1=15), matching style of O1
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For data rows 2-25: swap values in columns I & K, and in M & O,
# then fill new columns P and Q with 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # I
    $kVal = $ws.Cells.Item($r, 11).Value2  # K
    $mVal = $ws.Cells.Item($r, 13).Value2  # M
    $oVal = $ws.Cells.Item($r, 15).Value2  # O

    $ws.Cells.Item($r, 9).Value = $kVal    # I becomes old K
    $ws.Cells.Item($r, 11).Value = $iVal   # K becomes old I
    $ws.Cells.Item($r, 13).Value = $oVal   # M becomes old O
    $ws.Cells.Item($r, 15).Value = $mVal   # O becomes old M

    $ws.Cells.Item($r, 16).Value = 2       # P
    $ws.Cells.Item($r, 17).Value = 2       # Q
}

$wb.Save()
